$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: sub_emotion caring -> optimism, intensity moderate -> mild
$ws.Range("E2").Value = "optimism"
$ws.Range("F2").Value = "mild"

# Row 3: emotion fear -> happiness, sub_emotion nervousness -> excitement, intensity moderate -> mild
$ws.Range("D3").Value = "happiness"
$ws.Range("E3").Value = "excitement"
$ws.Range("F3").Value = "mild"

# Row 4: emotion happiness -> fear, sub_emotion optimism -> nervousness
$ws.Range("D4").Value = "fear"
$ws.Range("E4").Value = "nervousness"

# Row 5: emotion happiness -> fear, sub_emotion curiosity -> nervousness
$ws.Range("D5").Value = "fear"
$ws.Range("E5").Value = "nervousness"

# Row 6: intensity moderate -> mild
$ws.Range("F6").Value = "mild"
